$d = $word.ActiveDocument

# --- Part 1: merge split runs back into single runs. The source document had several
#     dates / phrases typed across multiple same-formatted runs (e.g. "Worked" + " " +
#     "with Multiple..."); the edit consolidates each into one run. Doing a Find &
#     Replace over the exact text reproduces Word's native run-coalescing behaviour. ---

$d.Content.Find.Execute("Worked with Multiple departments to coordinate and get all active product images fully updated for a storefront.", $false, $false, $false, $false, $false, $true, 1, $false, "Worked with Multiple departments to coordinate and get all active product images fully updated for a storefront.", 2) | Out-Null
$d.Content.Find.Execute("7-22-2019", $false, $false, $false, $false, $false, $true, 1, $false, "7-22-2019", 2) | Out-Null
$d.Content.Find.Execute("9-9-2019", $false, $false, $false, $false, $false, $true, 1, $false, "9-9-2019", 2) | Out-Null
$d.Content.Find.Execute("9-30-2019", $false, $false, $false, $false, $false, $true, 1, $false, "9-30-2019", 2) | Out-Null
$d.Content.Find.Execute("10-1-2019", $false, $false, $false, $false, $false, $true, 1, $false, "10-1-2019", 2) | Out-Null
$d.Content.Find.Execute("10-2-2019", $false, $false, $false, $false, $false, $true, 1, $false, "10-2-2019", 2) | Out-Null

# --- Part 2: insert new items 29-35 (each a number paragraph / Project paragraph /
#     Date paragraph, with assorted blank spacer paragraphs) immediately before the
#     trailing paragraph that carries the "_GoBack" bookmark. Inserting raw WordML at
#     the exact start of that paragraph's range -- with one extra empty <w:p/> tacked
#     onto the end of the payload to absorb into it -- adds the new paragraphs as
#     genuine siblings without disturbing the bookmark paragraph itself. ---

$count = $d.Paragraphs.Count
$pBookmark = $d.Paragraphs($count)
$startPos = $pBookmark.Range.Start
$insRange = $d.Range($startPos, $startPos)

$newItemsXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:wordDocument xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr><w:r><w:t>29</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Project: </w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>Worked on a site to update and test to ensure that ADA Compliance measures were met.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Date: </w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>1-11-2020</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr><w:r><w:t>30</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Project: </w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>Tested and worked with a team to see that server patching went smoothly.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Date: </w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>1-22-2020</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr><w:r><w:t>31</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Project: </w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>Worked to update cookie headers from the server side to prevent CSRF attacks.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Date: </w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>1-30-2020</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr><w:r><w:t>32</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Project: </w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t xml:space="preserve">Helped to diagnose and resolve an </w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>SQL Injection issue</w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t xml:space="preserve"> on a storefront</w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Date: </w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>1-31</w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>-2020</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr><w:r><w:t>33.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Project: </w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>Took over a process to update HR files for a site</w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Date: </w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>2-3</w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>-2020</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr><w:r><w:t>34</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Project: </w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>Setup contact lists for my department for remote working.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Date: </w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>3-11</w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>-2020</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr><w:r><w:t>35.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Project: </w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>QA new award emails before being sent across entire company</w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Date: </w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>4-14</w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>-2020</w:t></w:r></w:p><w:p></w:p></w:body></w:wordDocument>'
$insRange.InsertXML($newItemsXml)

# --- Part 3: append 9 trailing blank paragraphs after the bookmark paragraph. Since
#     it is the last paragraph in the body, "after" it means replacing its own full
#     range (through its end-of-paragraph mark) with [itself, unchanged] followed by
#     the new blank paragraphs. ---

$count = $d.Paragraphs.Count
$pBookmark = $d.Paragraphs($count)
$fullStart = $pBookmark.Range.Start
$fullEnd = $pBookmark.Range.End
$replRange = $d.Range($fullStart, $fullEnd)

$tailXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:wordDocument xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:hanging="720"/></w:pPr></w:p></w:body></w:wordDocument>'
$replRange.InsertXML($tailXml)

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
